$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell H1: "Save" — copy the formatting from the neighboring
# header cell (G1) so it reuses the existing bold/bordered header style,
# then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("H1").Value = "Save"

# New data column values for rows 2 and 3
$ws.Range("H2").Value = 1
$ws.Range("H3").Value = 1
